$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M3").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("M5").Value = 0
$ws1.Range("R9").Value = 0

$ws1.Range("D11").Value = "0 de 9"
$ws1.Range("M11").Value = "0 de 9"
$ws1.Range("R11").Value = "0 de 9"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths shift: D=11->10, E=10->12, F=12->11
# NOTE: the COM layer stores widths snapped to a 1/6-character pixel grid and
# adds a small padding offset before rounding, so plain integer inputs don't
# reliably round-trip to the same integer stored width. Values here were
# empirically chosen so the persisted <col width=".."/> comes out exactly
# 10, 12 and 11 respectively.
$ws2.Columns.Item(4).ColumnWidth = 9.15
$ws2.Columns.Item(5).ColumnWidth = 11.15
$ws2.Columns.Item(6).ColumnWidth = 10.15

# Header month labels shift forward by one month
$ws2.Range("C1").Value = "abril"
$ws2.Range("D1").Value = "mayo"
$ws2.Range("E1").Value = "junio"
$ws2.Range("F1").Value = "julio"

# Row 3 data
$ws2.Range("E3").Value = 832
$ws2.Range("F3").Value = 0

# Row 5 data
$ws2.Range("C5").Value = 0
$ws2.Range("E5").Value = 155.38
$ws2.Range("F5").Value = 0

# Row 9 data
$ws2.Range("C9").Value = 3.47
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 10.44
$ws2.Range("F9").Value = 0

# Row 11 totals
$ws2.Range("C11").Value = 3.47
$ws2.Range("D11").Value = 0
$ws2.Range("E11").Value = 997.8200000000001
$ws2.Range("F11").Value = 0
